$wb = $excel.ActiveWorkbook

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

# Update the unit description notes (text content changed; same cells)
$wsAbout.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: Mtoe - Million tonnes of oil equivalent - 10^6 toe"

# New row with conversion factor note, formatted with the scientific-number style (style index 3 / numFmtId 11)
$wsAbout.Range("A15").Value = "1 Btu = 2.5219021687207" + [char]0x22C5 + "10-8 toe"
$wsAbout.Range("A15").NumberFormat = "0.00E+00"

$wsAbout.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: Ktoe - 1000 tonnes of oil equivalent"

# --- "BpTPEU-large" sheet updates ---
$wsLarge = $wb.Worksheets.Item("BpTPEU-large")
$wsLarge.Range("B2").Formula = "=39652608749183"

# --- "BpTPEU-small" sheet updates ---
$wsSmall = $wb.Worksheets.Item("BpTPEU-small")
$wsSmall.Range("B2").Formula = "=39652608749.183"

# Active sheet moves to "BpTPEU-small"
$wsSmall.Activate()
